# Generate Report for Handback
# Update the handback-status timestamps to reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 0dad3cfe... row
$wsOverview.Range("G2").Value = "2016-09-01 15:31:38"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the 0dad3cfe... row
$wsZhCn.Range("H2").Value = "2016-09-01 15:31:33"
$wsZhCn.Range("K2").Value = "2016-09-01 15:32:16"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the 0dad3cfe... row
$wsDeDe.Range("H2").Value = "2016-09-01 15:31:38"
$wsDeDe.Range("K2").Value = "2016-09-01 15:32:29"
